$d = $word.ActiveDocument

# Update the nutrition table values (csv format fix):
# Cals: 2000 -> 2625
# Carbs: 52 -> 126
# Prot: -20 -> 39
# Fat: 22 -> 85

$table = $d.Tables.Item(1)

$table.Cell(2, 1).Range.Text = "2625"
$table.Cell(2, 2).Range.Text = "126"
$table.Cell(2, 3).Range.Text = "39"
$table.Cell(2, 4).Range.Text = "85"
